$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest crypto data.
# Price values are forced to literal text (leading apostrophe) because they are
# stored as plain strings in the source data (e.g. "29.476.93" is not a valid
# number) and assigning them as plain numbers would let Excel re-interpret /
# reformat values such as "8.910" or "245.43" (dropping trailing zeros or adding
# floating point noise).
$ws.Range("D2").Value = "'29.476.93"
$ws.Range("E2").Value = '  +2.14%  '
$ws.Range("D3").Value = "'1.854.70"
$ws.Range("E3").Value = '  +1.31%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = "'245.43"
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("D6").Value = "'0.6944"
$ws.Range("E6").Value = '  +0.72%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +0.49%  '
$ws.Range("D9").Value = "'0.07657"
$ws.Range("E9").Value = '  -0.53%  '
$ws.Range("D10").Value = "'23.51"
$ws.Range("E10").Value = '  +0.51%  '
$ws.Range("D11").Value = "'0.07758"
$ws.Range("E11").Value = '  -0.63%  '
$ws.Range("D12").Value = "'1.857.23"
$ws.Range("E12").Value = '  +1.71%  '
$ws.Range("D13").Value = "'5.139"
$ws.Range("E13").Value = '  +1.07%  '
$ws.Range("D14").Value = "'0.6934"
$ws.Range("E14").Value = '  +1.77%  '
$ws.Range("D15").Value = "'90.93"
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("D16").Value = "'6.306"
$ws.Range("E16").Value = '  -2.19%  '
$ws.Range("D17").Value = "'29.494.02"
$ws.Range("E17").Value = '  +2.23%  '
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").Value = "'2.102.78"
$ws.Range("E19").Value = '  +1.26%  '
$ws.Range("D20").Value = "'236.63"
$ws.Range("E20").Value = '  -2.13%  '
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").Value = "'7.624"
$ws.Range("E23").Value = '  +2.13%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("D25").Value = "'0.1484"
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = "'8.910"
$ws.Range("E26").Value = '  +1.39%  '
$ws.Range("D27").Value = "'159.33"
$ws.Range("E27").Value = '  +0.35%  '
$ws.Range("D28").Value = "'18.23"
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").Value = "'1.529"
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("D30").Value = "'4.247"
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("D31").Value = "'4.136"
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("E32").Value = '  +1.57%  '
$ws.Range("D33").Value = "'0.05242"
$ws.Range("E33").Value = '  +2.61%  '
$ws.Range("D34").Value = "'0.7774"
$ws.Range("E34").Value = '  +1.01%  '
$ws.Range("D35").Value = "'1.871"
$ws.Range("E35").Value = '  +1.16%  '
$ws.Range("E36").Value = '  +0.49%  '
$ws.Range("D37").Value = "'2.677"
$ws.Range("E37").Value = '  -0.50%  '
$ws.Range("D38").Value = "'1.317.75"
$ws.Range("E38").Value = '  +7.70%  '
$ws.Range("D39").Value = "'0.01869"
$ws.Range("E39").Value = '  +1.22%  '
$ws.Range("D40").Value = "'2.727"
$ws.Range("E40").Value = '  +1.29%  '
$ws.Range("D41").Value = "'0.9433"
$ws.Range("E41").Value = '  -1.33%  '
$ws.Range("D42").Value = "'106.17"
$ws.Range("E42").Value = '  -1.69%  '
$ws.Range("D43").Value = "'5.814"
$ws.Range("E43").Value = '  -0.31%  '
$ws.Range("E44").Value = '  +0.17%  '
$ws.Range("D45").Value = "'9.727"
$ws.Range("E45").Value = '  +0.88%  '
$ws.Range("D46").Value = "'2.002.53"
$ws.Range("E46").Value = '  +1.35%  '
$ws.Range("D47").Value = "'0.5232"
$ws.Range("E47").Value = '  +1.51%  '
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("E49").Value = '  +1.91%  '
$ws.Range("D50").Value = "'62.80"
$ws.Range("E50").Value = '  -2.44%  '
$ws.Range("D51").Value = "'0.05957"
$ws.Range("E51").Value = '  +1.05%  '
